# Day2 using Github clean stuff up
#
# Sheet "AB" (sheet1): add a new header column I ("WAEF - DEW") matching the
# existing header style, and add a new data row (row 3): EWEF / metering
# bridge / efw / AB / 54 / -110 / z / 2 / DW
#
# Sheet "BC" (sheet2): add a new data row (row 19): 243 / metering bridge /
# 2r43q / BC / 54 / -121 / 4r3qqqr3 / 3

$wb = $excel.ActiveWorkbook

$wsAB = $wb.Worksheets.Item("AB")
$wsBC = $wb.Worksheets.Item("BC")

# --- Sheet "AB": new header cell I2, copying the existing header formatting ---
$wsAB.Range("H2").Copy()
$wsAB.Range("I2").PasteSpecial(-4122)   # xlPasteFormats
$wsAB.Range("I2").Value = "WAEF - DEW"

# --- Sheet "AB": new data row 3 ---
$wsAB.Range("A3").Value = "EWEF"
$wsAB.Range("B3").Value = "metering bridge"
$wsAB.Range("C3").Value = "efw"
$wsAB.Range("D3").Value = "AB"
$wsAB.Range("E3").Value = 54
$wsAB.Range("F3").Value = -110
$wsAB.Range("G3").Value = "z"

# H3 ("2") must stay text, not become a number - force text format first
$wsAB.Range("H3").NumberFormat = "@"
$wsAB.Range("H3").Value = "2"
$wsAB.Range("H3").Style = "Normal"

$wsAB.Range("I3").Value = "DW"

# --- Sheet "BC": new data row 19 ---
# A19 ("243") must stay text, not become a number - force text format first
$wsBC.Range("A19").NumberFormat = "@"
$wsBC.Range("A19").Value = "243"
$wsBC.Range("A19").Style = "Normal"

$wsBC.Range("B19").Value = "metering bridge"
$wsBC.Range("C19").Value = "2r43q"
$wsBC.Range("D19").Value = "BC"
$wsBC.Range("E19").Value = 54
$wsBC.Range("F19").Value = -121
$wsBC.Range("G19").Value = "4r3qqqr3"

# H19 ("3") must stay text, not become a number - force text format first
$wsBC.Range("H19").NumberFormat = "@"
$wsBC.Range("H19").Value = "3"
$wsBC.Range("H19").Style = "Normal"
